# Added new login backup codes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sec-codes")

# The trailing block of codes (rows 13-16) is being replaced and moved up
# (the blank gap between the first and second block of codes shrinks), so
# clear those old rows out entirely first.
$ws.Range("A13:A16").ClearContents()

# First block of codes (rows 2-4) - replaced in place
$ws.Range("A2").Value = "1CZV4VTHGTVN"
$ws.Range("A3").Value = "51K0DF5KCN34"
$ws.Range("A4").Value = "HQT8HMXSF63S"

# Second block of codes, now starting at row 7 (was row 13), extended with
# two additional new codes at the end (rows 11-12)
$ws.Range("A7").Value = "5A41AVCX9PFR"
$ws.Range("A8").Value = "1YA1484DG5R7"
$ws.Range("A9").Value = "MJ6J3N01MN75"
$ws.Range("A10").Value = "B9AV6NE42R8W"
$ws.Range("A11").Value = "V2A2ZKV148W8"
$ws.Range("A12").Value = "DY9R3Z05BNS8"

# Highlight the header cell
$ws.Range("A1").Interior.ThemeColor = 10

# Update the saved cursor selection
[void]$ws.Range("B6").Select()
